# LOM3250.xlsx content restructuring.
#
# The row block A10:C24 is restructured: the "Objetivos:" value is
# replaced by the docente line, a new "Programa resumido: / Semestral"
# row is introduced, and everything below shifts up by one row, with
# the final (now-duplicate) trailing row removed so the sheet ends at
# row 23 instead of row 24.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to inject literal text that LOOKS like a date
# ("01/01/2012") without Excel's smart-typing turning it into a real
# date serial. A formula result pasted as values-only is taken
# verbatim as text.
$helper = "Z1"

# --- Row 10: Objetivos: / 519033 - Carlos Yujiro Shigue ---
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows(10).RowHeight = 60

# --- Row 11: Objectives: (no B/C value) ---
$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Rows(11).RowHeight = 60

# --- Row 12: Docentes responsáveis: (no B/C value) ---
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()

# --- Row 13: Programa resumido: / Semestral (new A cell) ---
$ws.Range("A10").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# --- Row 14: Short syllabus: (no B/C value) ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Rows(14).RowHeight = 60

# --- Row 15: Programa: / 01/01/2012 (new B/C cells, date-look text) ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B10").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range($helper).Formula = "=""01/01/2012"""
$ws.Range($helper).Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range($helper).ClearContents()
$excel.CutCopyMode = $false
$ws.Rows(15).RowHeight = 120

# --- Row 16: Syllabus: (no B/C value) ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Rows(16).RowHeight = 120

# --- Row 17: Avaliação: (no B/C value) ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

# --- Row 18: Método: / 519033 - Carlos Yujiro Shigue (new B/C cells) ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B10").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows(18).RowHeight = 60

# --- Row 19: Critério: / Em função da natureza... ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Em função da natureza deste curso, a avaliação será feita pela elaboração e apresentação de um plano de trabalho."
$ws.Range("C19").Value = "Em função da natureza deste curso, a avaliação será feita pela elaboração e apresentação de um plano de trabalho."
$ws.Rows(19).RowHeight = 60

# --- Row 20: Norma de recuperação: / Avaliação e atribuição... ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Avaliação e atribuição de nota do Trabalho de Graduação por uma comissão de professores."
$ws.Range("C20").Value = "Avaliação e atribuição de nota do Trabalho de Graduação por uma comissão de professores."
$ws.Rows(20).RowHeight = 60

# --- Row 21: Bibliografia: / A critério da banca... ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A critério da banca de avaliação poderá ser estabelecido um prazo para revisão e/ou correção da monografia."
$ws.Range("C21").Value = "A critério da banca de avaliação poderá ser estabelecido um prazo para revisão e/ou correção da monografia."
$ws.Rows(21).RowHeight = 120

# --- Row 22: Requisitos: (no B/C value) ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

# --- Row 23: (no A label) / LOM3238 requirement text (new B/C cells) ---
$ws.Range("A23").ClearContents()
$ws.Range("B13").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = "LOM3238 -  Projeto Integrado I  (Requisito)`n"
$ws.Range("C13").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "LOM3238 -  Projeto Integrado I  (Requisito)`n"
$ws.Rows(23).RowHeight = 30

$excel.CutCopyMode = $false

# --- Remove the old trailing row 24 so the used range shrinks to A1:C23 ---
$ws.Range("A24:C24").EntireRow.Delete()
